$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deduplicated nombre_aides (C) and montant_total (D) values per row,
# per the "[Fonds de solidarite] Deduplicate 2020-07-26 data" commit.
$data = @{
    2 = @(36399, 52682575)
    3 = @(88352, 129617155)
    4 = @(30275, 44864033)
    5 = @(8393, 12477676)
    6 = @(1845, 2743547)
    7 = @(138, 202093)
    11 = @(39783, 54060106)
    12 = @(9343, 13517559)
    13 = @(25310, 37138633)
    14 = @(8107, 12036168)
    15 = @(2084, 3101150)
    16 = @(390, 574123)
    19 = @(9842, 13070703)
    20 = @(13049, 18849366)
    21 = @(30972, 45481823)
    22 = @(10030, 14916231)
    23 = @(2558, 3807263)
    24 = @(469, 697592)
    25 = @(34, 50453)
    26 = @(11353, 15195673)
    27 = @(7407, 10731837)
    28 = @(21974, 32257013)
    29 = @(7610, 11325304)
    30 = @(1910, 2850544)
    31 = @(341, 508915)
    33 = @(8072, 10676850)
    34 = @(3066, 4421459)
    35 = @(7460, 10901189)
    36 = @(3008, 4456288)
    37 = @(791, 1179263)
    38 = @(144, 214232)
    40 = @(2305, 3112387)
    41 = @(16718, 24190852)
    42 = @(49734, 72943828)
    43 = @(18548, 27552530)
    44 = @(5420, 8072767)
    45 = @(1104, 1646792)
    49 = @(16176, 21572625)
    50 = @(1864, 2705440)
    51 = @(6450, 9489228)
    52 = @(2232, 3333324)
    53 = @(726, 1084305)
    54 = @(170, 251833)
    56 = @(6102, 8408761)
    57 = @(810, 1186584)
    58 = @(2018, 2994826)
    59 = @(823, 1226145)
    60 = @(281, 421258)
    61 = @(69, 103500)
    63 = @(1190, 1682894)
    64 = @(14924, 21566605)
    65 = @(43679, 63952099)
    66 = @(15355, 22828906)
    67 = @(4445, 6621560)
    68 = @(886, 1318596)
    71 = @(14679, 19387943)
    72 = @(48892, 71184838)
    73 = @(140072, 206454190)
    74 = @(61193, 91211890)
    75 = @(19467, 29090019)
    76 = @(4459, 6662472)
    83 = @(48504, 66201533)
    84 = @(4389, 6362228)
    85 = @(11120, 16342446)
    86 = @(3772, 5621415)
    87 = @(1312, 1960489)
    91 = @(5107, 6878553)
    92 = @(1487, 2149412)
    93 = @(4885, 7195584)
    94 = @(1864, 2777499)
    95 = @(659, 987460)
    99 = @(3302, 4383021)
    100 = @(565, 842964)
    101 = @(321, 479630)
    102 = @(113, 169500)
    105 = @(10451, 15182390)
    106 = @(28575, 41993886)
    107 = @(9573, 14236157)
    108 = @(2612, 3894307)
    109 = @(462, 690482)
    112 = @(9513, 12592630)
    113 = @(29479, 42543771)
    114 = @(64538, 94491545)
    115 = @(20884, 31047979)
    116 = @(5872, 8749280)
    117 = @(1067, 1595006)
    118 = @(68, 99420)
    121 = @(24990, 33422976)
    122 = @(34569, 49930096)
    123 = @(74367, 108823337)
    124 = @(23146, 34364029)
    125 = @(6166, 9164367)
    126 = @(1151, 1710911)
    130 = @(30480, 40548885)
    131 = @(12866, 18631818)
    132 = @(31590, 46421951)
    133 = @(11243, 16706059)
    134 = @(2880, 4294581)
    135 = @(463, 688490)
    138 = @(10515, 14044739)
    139 = @(33807, 48851883)
    140 = @(78820, 115526341)
    141 = @(23685, 35211145)
    142 = @(6177, 9220086)
    143 = @(1362, 2026192)
    146 = @(28222, 38171616)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item([int]$row, 3).Value = $vals[0]
    $ws.Cells.Item([int]$row, 4).Value = $vals[1]
}
